$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1790.5333
$ws.Range("I19").Value = 719.2
$ws.Range("J19").Value = 2326.2
$ws.Range("K19").Value = 719.2
$ws.Range("L19").Value = 2326.2
$ws.Range("M19").Value = -544.2
$ws.Range("N19").Value = -2676.2
$ws.Range("H40").Value = 2995.8
$ws.Range("J40").Value = 2993
$ws.Range("L40").Value = 2993
$ws.Range("N40").Value = -3343
$ws.Range("H58").Value = 1738.0769
$ws.Range("I58").Value = 363.375
$ws.Range("K58").Value = 1090.125
$ws.Range("M58").Value = -940.125
$ws.Range("H86").Value = 713.2
$ws.Range("I86").Value = 713.2
$ws.Range("K86").Value = 713.2
$ws.Range("M86").Value = 409.8
$ws.Range("H89").Value = 713.2
$ws.Range("I89").Value = 713.2
$ws.Range("K89").Value = 3566
$ws.Range("M89").Value = 2050
$ws.Range("H100").Value = 4094.75
$ws.Range("I100").Value = 400
$ws.Range("J100").Value = 5326.3335
$ws.Range("K100").Value = 400
$ws.Range("L100").Value = 5326.3335
$ws.Range("M100").Value = 141
$ws.Range("N100").Value = -6408.3335
$ws.Range("H116").Value = 8306.799999999999
$ws.Range("I116").Value = 16311.571
$ws.Range("J116").Value = 3996.5386
$ws.Range("K116").Value = 16311.571
$ws.Range("L116").Value = 3996.5386
$ws.Range("M116").Value = -12869.571
$ws.Range("N116").Value = -10880.5386
$ws.Range("H137").Value = 30897.824
$ws.Range("I137").Value = 1001.3077
$ws.Range("J137").Value = 49405.19
$ws.Range("K137").Value = 3003.9231
$ws.Range("L137").Value = 148215.57
$ws.Range("M137").Value = -453.9231
$ws.Range("N137").Value = -153315.57
$ws.Range("H138").Value = 3643.6035
$ws.Range("I138").Value = 3692.1667
$ws.Range("J138").Value = 3621.75
$ws.Range("K138").Value = 11076.5001
$ws.Range("L138").Value = 10865.25
$ws.Range("M138").Value = -5936.500100000001
$ws.Range("N138").Value = -21145.25

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2951.4675
$ws.Range("I32").Value = 2347.9092
$ws.Range("J32").Value = 6572.8184
$ws.Range("K32").Value = 2347.9092
$ws.Range("L32").Value = 6572.8184
$ws.Range("M32").Value = -2060.9092
$ws.Range("N32").Value = -7146.8184
$ws.Range("H61").Value = 2092.6365
$ws.Range("I61").Value = 1077.25
$ws.Range("J61").Value = 4800.3335
$ws.Range("K61").Value = 1077.25
$ws.Range("L61").Value = 4800.3335
$ws.Range("M61").Value = -865.25
$ws.Range("N61").Value = -5224.3335
$ws.Range("H74").Value = 937.1818
$ws.Range("I74").Value = 942.1429000000001
$ws.Range("J74").Value = 928.5
$ws.Range("K74").Value = 942.1429000000001
$ws.Range("L74").Value = 928.5
$ws.Range("M74").Value = -68.14290000000005
$ws.Range("N74").Value = -2676.5
$ws.Range("H77").Value = 937.1818
$ws.Range("I77").Value = 942.1429000000001
$ws.Range("J77").Value = 928.5
$ws.Range("K77").Value = 4710.7145
$ws.Range("L77").Value = 4642.5
$ws.Range("M77").Value = -342.7145
$ws.Range("N77").Value = -13378.5
$ws.Range("H102").Value = 1999.6666
$ws.Range("I102").Value = 1999.6666
$ws.Range("K102").Value = 1999.6666
$ws.Range("M102").Value = -377.6666
$ws.Range("H110").Value = 1600.8889
$ws.Range("I110").Value = 1058.2858
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 1058.2858
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = 986.7141999999999
$ws.Range("N110").Value = -7590
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H132").Value = 1932.6274
$ws.Range("I132").Value = 1462.3928
$ws.Range("J132").Value = 2505.087
$ws.Range("K132").Value = 4387.178400000001
$ws.Range("L132").Value = 7515.261
$ws.Range("M132").Value = -1857.178400000001
$ws.Range("N132").Value = -12575.261
$ws.Range("H136").Value = 2092.6365
$ws.Range("I136").Value = 1077.25
$ws.Range("J136").Value = 4800.3335
$ws.Range("K136").Value = 3231.75
$ws.Range("L136").Value = 14401.0005
$ws.Range("M136").Value = -681.75
$ws.Range("N136").Value = -19501.0005

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1892.1666
$ws.Range("I105").Value = 2009.1333
$ws.Range("J105").Value = 1307.3334
$ws.Range("K105").Value = 2009.1333
$ws.Range("L105").Value = 1307.3334
$ws.Range("M105").Value = -262.1333
$ws.Range("N105").Value = -4801.3334
$ws.Range("H134").Value = 5127.478
$ws.Range("I134").Value = 5350.359
$ws.Range("J134").Value = 3885.7144
$ws.Range("K134").Value = 16051.077
$ws.Range("L134").Value = 11657.1432
$ws.Range("M134").Value = -13516.077
$ws.Range("N134").Value = -16727.1432

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2116.158
$ws.Range("I31").Value = 1015.0769
$ws.Range("J31").Value = 4501.8335
$ws.Range("K31").Value = 1015.0769
$ws.Range("L31").Value = 4501.8335
$ws.Range("M31").Value = -720.0769
$ws.Range("N31").Value = -5091.8335
$ws.Range("H34").Value = 2116.158
$ws.Range("I34").Value = 1015.0769
$ws.Range("J34").Value = 4501.8335
$ws.Range("K34").Value = 1015.0769
$ws.Range("L34").Value = 4501.8335
$ws.Range("M34").Value = -813.0769
$ws.Range("N34").Value = -4905.8335
$ws.Range("H58").Value = 1319208.1
$ws.Range("I58").Value = 1891008.5
$ws.Range("J58").Value = 4067.4
$ws.Range("K58").Value = 1891008.5
$ws.Range("L58").Value = 4067.4
$ws.Range("M58").Value = -1890805.5
$ws.Range("N58").Value = -4473.4
$ws.Range("H99").Value = 1601.4445
$ws.Range("I99").Value = 1334
$ws.Range("K99").Value = 1334
$ws.Range("M99").Value = 164
$ws.Range("H105").Value = 1998
$ws.Range("I105").Value = 1998
$ws.Range("K105").Value = 1998
$ws.Range("M105").Value = -251
$ws.Range("H126").Value = 1601.4445
$ws.Range("I126").Value = 1334
$ws.Range("K126").Value = 4002
$ws.Range("M126").Value = -1532
$ws.Range("H132").Value = 2769.8147
$ws.Range("I132").Value = 984.86664
$ws.Range("J132").Value = 5001
$ws.Range("K132").Value = 2954.59992
$ws.Range("L132").Value = 15003
$ws.Range("M132").Value = -424.5999199999997
$ws.Range("N132").Value = -20063
$ws.Range("H134").Value = 2096.4688
$ws.Range("I134").Value = 1340.7307
$ws.Range("J134").Value = 5371.3335
$ws.Range("K134").Value = 4022.1921
$ws.Range("L134").Value = 16114.0005
$ws.Range("M134").Value = -1487.1921
$ws.Range("N134").Value = -21184.0005
$ws.Range("H136").Value = 1319208.1
$ws.Range("I136").Value = 1891008.5
$ws.Range("J136").Value = 4067.4
$ws.Range("K136").Value = 5673025.5
$ws.Range("L136").Value = 12202.2
$ws.Range("M136").Value = -5670475.5
$ws.Range("N136").Value = -17302.2

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 328.7
$ws.Range("I5").Value = 301.27274
$ws.Range("J5").Value = 344.57895
$ws.Range("K5").Value = 903.81822
$ws.Range("L5").Value = 1033.73685
$ws.Range("M5").Value = -791.81822
$ws.Range("N5").Value = -1257.73685
$ws.Range("H107").Value = 2679.76
$ws.Range("J107").Value = 2731.7317
$ws.Range("L107").Value = 8195.195099999999
$ws.Range("N107").Value = -12035.1951
$ws.Range("H131").Value = 14306959
$ws.Range("I131").Value = 38462024
$ws.Range("J131").Value = 33512.047
$ws.Range("K131").Value = 115386072
$ws.Range("L131").Value = 100536.141
$ws.Range("M131").Value = -115381032
$ws.Range("N131").Value = -110616.141
$ws.Range("H132").Value = 1652.3334
$ws.Range("I132").Value = 763.3333
$ws.Range("J132").Value = 1874.5834
$ws.Range("K132").Value = 6869.9997
$ws.Range("L132").Value = 16871.2506
$ws.Range("M132").Value = -4339.9997
$ws.Range("N132").Value = -21931.2506
$ws.Range("H133").Value = 31254124
$ws.Range("I133").Value = 125001500
$ws.Range("K133").Value = 375004500
$ws.Range("M133").Value = -374999440
$ws.Range("H135").Value = 328.7
$ws.Range("I135").Value = 301.27274
$ws.Range("J135").Value = 344.57895
$ws.Range("K135").Value = 2711.45466
$ws.Range("L135").Value = 3101.21055
$ws.Range("M135").Value = -176.4546599999999
$ws.Range("N135").Value = -8171.21055

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1412.5
$ws.Range("I113").Value = 1320
$ws.Range("K113").Value = 1320
$ws.Range("M113").Value = 850
$ws.Range("H132").Value = 1542338.8
$ws.Range("I132").Value = 4809601
$ws.Range("J132").Value = 4803.7646
$ws.Range("K132").Value = 14428803
$ws.Range("L132").Value = 14411.2938
$ws.Range("M132").Value = -14426273
$ws.Range("N132").Value = -19471.2938

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9324.200000000001
$ws.Range("I16").Value = 10780.25
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 10780.25
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = -10610.25
$ws.Range("N16").Value = -3840
$ws.Range("H46").Value = 1156.4667
$ws.Range("J46").Value = 1303.9166
$ws.Range("L46").Value = 1303.9166
$ws.Range("N46").Value = -1679.9166
$ws.Range("H93").Value = 784.05
$ws.Range("I93").Value = 805.3889
$ws.Range("J93").Value = 592
$ws.Range("K93").Value = 805.3889
$ws.Range("L93").Value = 592
$ws.Range("M93").Value = 442.6111
$ws.Range("N93").Value = -3088
$ws.Range("H122").Value = 4870.52
$ws.Range("I122").Value = 3397.8235
$ws.Range("K122").Value = 10193.4705
$ws.Range("M122").Value = -7743.470499999999
$ws.Range("H132").Value = 2701.3408
$ws.Range("I132").Value = 1043.4166
$ws.Range("J132").Value = 4690.85
$ws.Range("K132").Value = 3130.2498
$ws.Range("L132").Value = 14072.55
$ws.Range("M132").Value = -600.2498000000001
$ws.Range("N132").Value = -19132.55
$ws.Range("H136").Value = 3883.1714
$ws.Range("I136").Value = 2377.35
$ws.Range("J136").Value = 5890.933
$ws.Range("K136").Value = 7132.049999999999
$ws.Range("L136").Value = 17672.799
$ws.Range("M136").Value = -4582.049999999999
$ws.Range("N136").Value = -22772.799

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1186.6595
$ws.Range("I132").Value = 466.5
$ws.Range("J132").Value = 2247.9473
$ws.Range("K132").Value = 1399.5
$ws.Range("L132").Value = 6743.841899999999
$ws.Range("M132").Value = 1130.5
$ws.Range("N132").Value = -11803.8419
$ws.Range("H136").Value = 12080014
$ws.Range("I136").Value = 23150782
$ws.Range("J136").Value = 2811.7727
$ws.Range("K136").Value = 69452346
$ws.Range("L136").Value = 8435.3181
$ws.Range("M136").Value = -69449796
$ws.Range("N136").Value = -13535.3181
